$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value set
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting all subsequent rows up by one.
$ws.Rows.Item(11).Delete()

# After the shift, the "Case Sensitive" row (now row 14) gets a text value of "true"
# (leading apostrophe forces this to be stored as text, not a boolean)
$ws.Range("B14").Value = "'true"
